$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.713.16"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "3.701.28"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("D4").Value = "'0.999"

$ws.Range("D5").Value = "'672.04"

$ws.Range("D6").Value = "'161.87"
$ws.Range("E6").Value = "  +1.61%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("D12").Value = "'0.0000236"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("D14").Value = "3.726.52"
$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Value = "69.690.87"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("E16").Value = "  +1.62%  "

$ws.Range("D17").Value = "'16.24"
$ws.Range("E17").Value = "  +2.53%  "

$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").Value = "'474.98"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").Value = "'80.45"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").Value = "3.848.68"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("E24").Value = "  +3.11%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.166"
$ws.Range("E33").Value = "  +3.88%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.92"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").Value = "3.689.41"
$ws.Range("E35").Value = "  +1.05%  "

$ws.Range("E36").Value = "  +4.71%  "

$ws.Range("D37").Value = "'6.11"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D39").Value = "'2.25"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'0.0910"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").Value = "'172.37"
$ws.Range("E42").Value = "  +3.48%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("E45").Value = "  +2.68%  "

$ws.Range("D46").Value = "'0.000282"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  +0.78%  "
